$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.514003396034241
$ws.Range("B1").Value = 6.474394798278809
$ws.Range("C1").Value = 3.475831747055054
$ws.Range("D1").Value = 1.550420522689819
$ws.Range("E1").Value = 1.091705560684204
